$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Column D holds text-formatted price strings (e.g. "1.00", "69.042.07").
# Force text number-format before assigning so Excel does not silently
# coerce numeric-looking strings ("1.00", "680.35", ...) into numbers.
foreach ($addr in @('D2','D3','D4','D5','D6','D7','D9','D11','D14','D15','D16','D17','D19','D21','D22','D23','D24','D25','D26','D27','D28','D29','D31','D33','D36','D37','D39','D40','D42','D44','D45','D46','D48','D49','D51')) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '69.042.07'
$ws.Range('E2').Value = '  -2.20%  '
$ws.Range('D3').Value = '3.675.89'
$ws.Range('E3').Value = '  -3.01%  '
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').Value = '680.35'
$ws.Range('E5').Value = '  -4.13%  '
$ws.Range('D6').Value = '161.94'
$ws.Range('E6').Value = '  -4.65%  '
$ws.Range('D7').Value = '3.675.58'
$ws.Range('E7').Value = '  -2.99%  '
$ws.Range('E8').Value = '  +0.18%  '
$ws.Range('D9').Value = '0.495'
$ws.Range('E9').Value = '  -4.99%  '
$ws.Range('E10').Value = '  -7.54%  '
$ws.Range('D11').Value = '7.23'
$ws.Range('E11').Value = '  -1.55%  '
$ws.Range('E12').Value = '  -1.40%  '
$ws.Range('E13').Value = '  -7.18%  '
$ws.Range('D14').Value = '33.34'
$ws.Range('E14').Value = '  -7.77%  '
$ws.Range('D15').Value = '4.298.45'
$ws.Range('E15').Value = '  -2.96%  '
$ws.Range('D16').Value = '3.682.61'
$ws.Range('E16').Value = '  -2.98%  '
$ws.Range('D17').Value = '69.151.45'
$ws.Range('E17').Value = '  -2.08%  '
$ws.Range('E18').Value = '  -1.79%  '
$ws.Range('D19').Value = '16.26'
$ws.Range('E19').Value = '  -6.28%  '
$ws.Range('E20').Value = '  -7.34%  '
$ws.Range('D21').Value = '482.22'
$ws.Range('E21').Value = '  -2.11%  '
$ws.Range('D22').Value = '9.78'
$ws.Range('E22').Value = '  -7.70%  '
$ws.Range('D23').Value = '0.662'
$ws.Range('E23').Value = '  -8.72%  '
$ws.Range('D24').Value = '79.46'
$ws.Range('E24').Value = '  -6.33%  '
$ws.Range('D25').Value = '3.822.31'
$ws.Range('E25').Value = '  -2.99%  '
$ws.Range('B26').Value = 'InternetComputer(DFINITY)'
$ws.Range('C26').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D26').Value = '11.53'
$ws.Range('E26').Value = '  -4.39%  '
$ws.Range('B27').Value = 'Dai'
$ws.Range('C27').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D27').Value = '0.999'
$ws.Range('E27').Value = '  -0.09%  '
$ws.Range('B28').Value = 'PEPE'
$ws.Range('C28').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D28').Value = '0.0000127'
$ws.Range('E28').Value = '  -12.20%  '
$ws.Range('D29').Value = '9.49'
$ws.Range('E29').Value = '  -9.29%  '
$ws.Range('E30').Value = '  -10.81%  '
$ws.Range('D31').Value = '2.74'
$ws.Range('E31').Value = '  -11.57%  '
$ws.Range('E32').Value = '  -5.17%  '
$ws.Range('D33').Value = '6.73'
$ws.Range('E33').Value = '  -7.99%  '
$ws.Range('E34').Value = '  +0.18%  '
$ws.Range('D36').Value = '0.163'
$ws.Range('E36').Value = '  -6.59%  '
$ws.Range('D37').Value = '3.644.27'
$ws.Range('E37').Value = '  -3.03%  '
$ws.Range('E38').Value = '  -5.98%  '
$ws.Range('D39').Value = '6.04'
$ws.Range('E39').Value = '  +2.21%  '
$ws.Range('D40').Value = '0.0939'
$ws.Range('E40').Value = '  -7.15%  '
$ws.Range('D42').Value = '2.17'
$ws.Range('E42').Value = '  -6.33%  '
$ws.Range('E43').Value = '  +0.09%  '
$ws.Range('D44').Value = '0.955'
$ws.Range('E44').Value = '  -8.00%  '
$ws.Range('D45').Value = '156.75'
$ws.Range('E45').Value = '  -4.68%  '
$ws.Range('D46').Value = '47.98'
$ws.Range('E46').Value = '  -1.88%  '
$ws.Range('E47').Value = '  -14.92%  '
$ws.Range('B48').Value = 'Bittensor'
$ws.Range('C48').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D48').Value = '391.52'
$ws.Range('E48').Value = '  -6.65%  '
$ws.Range('B49').Value = 'FLOKI'
$ws.Range('C49').Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
$ws.Range('D49').Value = '0.000276'
$ws.Range('E49').Value = '  -11.30%  '
$ws.Range('E50').Value = '  -4.94%  '
$ws.Range('D51').Value = '28.25'
$ws.Range('E51').Value = '  +1.52%  '
